# Insert a new data row into the "Vega Modelo de Temuco - Coliflor" sheet.
# The new row is inserted at row 366, pushing the existing rows 366-476 down
# to 367-477 (Excel's Insert() shifts cells down and copies formatting from
# the row above, matching the row template used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 366.
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with its data.
$ws.Cells.Item(366, 1).Value = 10
$ws.Cells.Item(366, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(366, 3).Value = "La Araucanía"
$ws.Cells.Item(366, 4).Value = 44841
$ws.Cells.Item(366, 5).Value = 9
$ws.Cells.Item(366, 6).Value = 100112008
$ws.Cells.Item(366, 7).Value = "Coliflor"
$ws.Cells.Item(366, 8).Value = "Sin especificar"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 2000
$ws.Cells.Item(366, 11).Value = 1200
$ws.Cells.Item(366, 12).Value = 1200
$ws.Cells.Item(366, 13).Value = 1200
$ws.Cells.Item(366, 14).Value = "$/unidad"
$ws.Cells.Item(366, 15).Value = "Región Metropolitana"
$ws.Cells.Item(366, 16).Value = 1200
$ws.Cells.Item(366, 17).Value = 1
$ws.Cells.Item(366, 18).Value = "Hortaliza"
